# Update Num_Inclusions (column C) values to use the 3rd quartile instead
# of the mean for the underlying per-cell inclusion counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = 8
    4   = 4
    6   = 6
    12  = 11
    15  = 11
    17  = 0
    18  = 4
    32  = 0
    37  = 2
    40  = 3
    57  = 9
    58  = 2
    59  = 1
    60  = 1
    73  = 1
    101 = 4
    108 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
